$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every changed cell originally held an inline-string (text) value, even
# though several look like plain numbers (e.g. "94.35") or percentages.
# A bare `.Value = '94.35'` gets auto-coerced to a numeric cell by the
# engine, so a leading apostrophe forces text entry (standard Excel "treat
# as text" convention); the apostrophe itself is not stored in the value.
# The Style reset afterwards drops the resulting quote-prefix formatting so
# the cell keeps the default (unstyled) appearance, matching the source file.

$ws.Range('D2').Value = "'43.852.91"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = "'  +0.46%  "
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = "'2.236.37"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = "'  +1.79%  "
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = "'  +0.02%  "
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = "'270.88"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "'  +4.20%  "
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = "'94.35"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = "'  +14.83%  "
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').Value = "'0.623"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = "'  +1.04%  "
$ws.Range('E7').Style = "Normal"
$ws.Range('E8').Value = "'  -0.03%  "
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').Value = "'0.642"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = "'  +8.14%  "
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = "'46.40"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = "'  +6.59%  "
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = "'0.0956"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = "'  +4.03%  "
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').Value = "'8.38"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = "'  +20.87%  "
$ws.Range('E12').Style = "Normal"
$ws.Range('E13').Value = "'  +2.03%  "
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = "'15.33"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "'  +7.43%  "
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = "'2.571.58"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = "'  +1.84%  "
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = "'0.824"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "'  +6.11%  "
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = "'2.238.22"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = "'  +1.48%  "
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = "'43.830.17"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = "'  +0.56%  "
$ws.Range('E18').Style = "Normal"
$ws.Range('E19').Value = "'  +1.78%  "
$ws.Range('E19').Style = "Normal"
$ws.Range('E20').Value = "'  +4.48%  "
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = "'70.84"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = "'  +1.46%  "
$ws.Range('E21').Style = "Normal"
$ws.Range('E22').Value = "'  -4.24%  "
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = "'234.61"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = "'  +1.64%  "
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = "'9.16"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = "'  +3.30%  "
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = "'0.999"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "'  -0.05%  "
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = "'11.40"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = "'  +6.38%  "
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = "'2.51"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = "'  +12.18%  "
$ws.Range('E27').Style = "Normal"
$ws.Range('E28').Value = "'  +6.22%  "
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = "'40.43"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = "'  -4.65%  "
$ws.Range('E29').Style = "Normal"
$ws.Range('E30').Value = "'  +2.85%  "
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = "'172.68"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = "'  -0.73%  "
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = "'0.0911"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = "'  +3.91%  "
$ws.Range('E32').Style = "Normal"
$ws.Range('E33').Value = "'  +2.81%  "
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = "'5.51"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = "'  +3.43%  "
$ws.Range('E34').Style = "Normal"
$ws.Range('E35').Value = "'  +1.98%  "
$ws.Range('E35').Style = "Normal"
$ws.Range('E36').Value = "'  -0.81%  "
$ws.Range('E36').Style = "Normal"
$ws.Range('E37').Value = "'  +0.19%  "
$ws.Range('E37').Style = "Normal"
$ws.Range('E38').Value = "'  -3.66%  "
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').Value = "'3.56"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = "'  +23.78%  "
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = "'12.82"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = "'  -2.13%  "
$ws.Range('E40').Style = "Normal"
$ws.Range('E41').Value = "'  +12.52%  "
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = "'2.15"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = "'  +2.32%  "
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').Value = "'63.44"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = "'  -1.04%  "
$ws.Range('E43').Style = "Normal"
$ws.Range('E44').Value = "'  -1.00%  "
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = "'0.0996"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = "'  +1.84%  "
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').Value = "'101.18"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = "'  +0.77%  "
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = "'8.41"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = "'  +1.39%  "
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').Value = "'1.15"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = "'  +3.97%  "
$ws.Range('E48').Style = "Normal"
$ws.Range('E49').Value = "'  +2.67%  "
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = "'0.451"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "'  +2.54%  "
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = "'2.456.10"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = "'  +1.82%  "
$ws.Range('E51').Style = "Normal"
